$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Save" in H1, copying the format from G1 (bold/bordered header style)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the "Save" indicator values for rows 2-7
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 0
